# The workbook's A2 cell held the instance id " i-0825f73d564da4ff4"
# with a stray leading space. Re-enter the value without the leading
# space (Excel dedupes/re-appends the shared string, which is why the
# shared-strings table gets reordered on save) and leave the selection
# on A2, matching the author's final click position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "i-0825f73d564da4ff4"

$ws.Range("A2").Select()
